# Updates the cryptos price/coin list to reflect the latest scrape.
# Column D values are numeric-looking strings that must remain stored
# as text (matching the source data which stores everything as text),
# so NumberFormat "@" is applied before assignment and the style is
# reset back to Normal afterwards so no stray cell style is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '264.25'
$ws.Range('D2').Style = "Normal"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '22.78'
$ws.Range('D3').Style = "Normal"
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '6.211'
$ws.Range('D4').Style = "Normal"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.06119'
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '3.541'
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '6.733'
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.372'
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.8169'
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1593'
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08210'
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.03367'
$ws.Range('D12').Style = "Normal"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.03154'
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.09275'
$ws.Range('D14').Style = "Normal"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.904'
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.001720'
$ws.Range('D16').Style = "Normal"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.04840'
$ws.Range('D17').Style = "Normal"
$ws.Range('B18').Value = 'TigerCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.006242'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '17TigerCashTCH'
$ws.Range('B19').Value = 'BitKan'
$ws.Range('C19').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.001102'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '18BitKanKAN'
$ws.Range('B20').Value = 'HotbitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.003203'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '19HotbitTokenHTBWorstin24h'
$ws.Range('B21').Value = 'NitroEx'
$ws.Range('C21').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0001504'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '20NitroExNTX'
$ws.Range('B22').Value = 'LEO'
$ws.Range('C22').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.696'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '21LEOLEO'
$ws.Range('B23').Value = 'BTSEToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.265'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '22BTSETokenBTSE'
$ws.Range('B24').Value = 'One'
$ws.Range('C24').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.01343'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '23OneONE'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.3407'
$ws.Range('D25').Style = "Normal"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0002688'
$ws.Range('D27').Style = "Normal"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.04654'
$ws.Range('D40').Style = "Normal"
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.007269'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.003910'
$ws.Range('D42').Style = "Normal"
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1124'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '42BKEXTokenBKK'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.01036'
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00006177'
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.00000000752'
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.7520'
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.1679'
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.00002106'
$ws.Range('D49').Style = "Normal"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.01243'
$ws.Range('D50').Style = "Normal"
